# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (value
# "stock" for every data row), inserted right after the existing "total"
# column and before the "date" column. Inserting a real column shifts the
# old date/legislator_name/legislator_id columns one slot to the right and
# carries their formatting/values along automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Previous layout: B=name C=owner D=quantity E=face_value F=currency
# G=total H=date I=legislator_name J=legislator_id
# Insert a new column at H so the old H (date) becomes I, etc.
$ws.Columns.Item(8).Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
